# Danger zone game flow, some fixes
# - Insert a new "food_beefStew" row before the TCS food description rows
# - Append three new "dangerZoneTrivia" rows at the bottom of the sheet
# - Update the sheet view's selection to the new last cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 37 (old rows 37-67 shift down to 38-68)
$ws.Rows("37:37").Insert()

# Populate new/changed cells in the exact order the strings were originally
# authored so the shared-string table indices line up with the source file.
$ws.Cells.Item(69, 1).Value = "dangerZoneTrivia1"
$ws.Cells.Item(37, 1).Value = "food_beefStew"
$ws.Cells.Item(37, 2).Value = "Beef Stew"
$ws.Cells.Item(70, 1).Value = "dangerZoneTrivia10"
$ws.Cells.Item(71, 1).Value = "dangerZoneTrivia11"
$ws.Cells.Item(69, 2).Value = "A pot of beef stew has been placed in an ice-water bath to cool. After constant stirring, three hours have passed, and the temperature reads 70° F. Is the beef stew ready to be placed into the cooler?"
$ws.Cells.Item(71, 2).Value = "Yes, the stew has reached the proper temperature of 70° F, it should be ready to be placed in the cooler."
$ws.Cells.Item(70, 2).Value = "No, the temperature didn't reach 70° F before two hours, it should be thrown out."

# Update the view's selection to match the new bottom row of data
[void]$ws.Range("B71").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 55
